# Generate Report for Handback
#
# - "Status" column value (shared across Overview!E3/F3, zh-cn!C3, de-de!C3)
#   changes from "Ready for handoff" to "Handback transform failed".
# - "Error Detail" column (P) on the zh-cn and de-de sheets gets a new
#   diagnostic message explaining the handback/handoff file name mismatch.
# - The "Error Detail" column is widened to fit the new long message.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "Ready for handoff" -> "Handback transform failed" -----------
$newStatus = "Handback transform failed"
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Error Detail messages for the failed handback transform --------------
$zhCnError = "Handback file name: 1eyieywm.kp1 is different with handoff file name: 6e1ca166-6b50-4ff7-89eb-657d67708ed9.d52d505bc12fa8e57a697c7f97197e7c5a5a9023.zh-cn."
$deDeError = "Handback file name: 1eyieywm.kp1 is different with handoff file name: 6e1ca166-6b50-4ff7-89eb-657d67708ed9.d52d505bc12fa8e57a697c7f97197e7c5a5a9023.de-de."

$wsZhCn.Range("P3").Value = $zhCnError
$wsDeDe.Range("P3").Value = $deDeError

# --- Widen the Error Detail column (P) to fit the new messages ------------
$fullWidth = $wsZhCn.Columns.Item(1).ColumnWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $fullWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $fullWidth
